$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STM32F429VE")

# Commit: use spi4 soft cs -> pin PE4 (SPI4_NSS) silkscreen label set to "软cs"
$ws.Range("E7").Value = "软cs"

# Update active selection as seen in the diff (cosmetic, matches target view state)
$ws.Range("E9").Select()
